$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Range("E11").NumberFormat = "#,##0"
Write-Host "done"
